# Wrap the Excel write in try/except (per commit message "added try except
# to excel write function") and refresh each regression-summary sheet's
# timestamp from the original run ("Sun, 29 Dec 2019 16:11:22/23") to the
# re-run timestamp ("Wed, 01 Jan 2020 23:19:00").

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    try {
        $text = $ws.Range("B2").Text
        if ($text -and $text.Contains("Date:")) {
            $origRowHeight = $ws.Rows.Item(2).RowHeight
            $updated = $text.Replace("Sun, 29 Dec 2019", "Wed, 01 Jan 2020")
            $updated = $updated.Replace("16:11:22", "23:19:00")
            $updated = $updated.Replace("16:11:23", "23:19:00")
            $ws.Range("B2").Value = $updated
            # writing the (same-length) text shouldn't reflow the row, but
            # restore the original wrapped-text row height just in case
            $ws.Rows.Item(2).RowHeight = $origRowHeight
        }
    } catch {
        Write-Host "Failed to update sheet $($ws.Name): $_"
    }
}
